# Update "想去人数" (F column) figures for several events.
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the
# same events, so both need to be updated in lockstep.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 81
$ws1.Range("F4").Value = 253
$ws1.Range("F6").Value = 238
$ws1.Range("F8").Value = 1893
$ws1.Range("F9").Value = 340
$ws1.Range("F10").Value = 4394
$ws1.Range("F11").Value = 63
$ws1.Range("F12").Value = 310

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 81
$ws4.Range("F6").Value = 253
$ws4.Range("F8").Value = 238
$ws4.Range("F12").Value = 1893
$ws4.Range("F13").Value = 340
$ws4.Range("F14").Value = 4394
$ws4.Range("F15").Value = 63
$ws4.Range("F16").Value = 310
